$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 67

# Columns A (date-like) and C (numeric-like) need a leading apostrophe so
# Excel stores the literal text instead of coercing to a date/number.
$ws.Cells.Item($row, 1).Value = "'2025-11-22"
$ws.Cells.Item($row, 2).Value = "Pick 3"
$ws.Cells.Item($row, 3).Value = "'251122"
$ws.Cells.Item($row, 4).Value = "0-2-1"
$ws.Cells.Item($row, 5).Value = "2025-11-22T21:36:33.432+04:00"
